# Fixed issue #13 Permitir que en los ficheros de metadatos dos columnas
# se puedan relacionar para crear SKOS jerarquicos.
#
# A new "machine-readable slug" row is inserted right after the header
# row (row 1), containing a slugified version of each column header
# (lowercase, hyphen-separated, no accents). All the rows that used to
# follow the header (the sdmx-dimension / dim / URI metadata rows) are
# pushed down one position. The old trailing row 5, which only held a
# stray "mapping-ano.xlsx" value in column K, is removed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 2, pushing existing rows 2-5 down
# to rows 3-6 (this keeps their formatting/style intact).
$ws.Rows.Item(2).Insert()

# Populate the newly inserted row 2 with the machine-readable slug for
# each column, matching the header text in row 1.
$ws.Cells.Item(2, 1).Value  = "ccaa-nombre"
$ws.Cells.Item(2, 2).Value  = "comarca-nombre"
$ws.Cells.Item(2, 3).Value  = "comarca-codigo"
$ws.Cells.Item(2, 4).Value  = "tipo-licencias-descripcion"
$ws.Cells.Item(2, 5).Value  = "ccaa-codigo"
$ws.Cells.Item(2, 6).Value  = "tipo-licencias-codigo"
$ws.Cells.Item(2, 7).Value  = "licencias"
$ws.Cells.Item(2, 8).Value  = "provincia-codigo"
$ws.Cells.Item(2, 9).Value  = "municipio-codigo"
$ws.Cells.Item(2, 10).Value = "provincia-nombre"
$ws.Cells.Item(2, 11).Value = "ano"
$ws.Cells.Item(2, 12).Value = "municipio-nombre"

# The old row 5 (now shifted to row 6) only contained a leftover
# "mapping-ano.xlsx" value in column K and is no longer part of the
# table; remove it entirely.
$ws.Rows.Item(6).Delete()
